# aggiornamento fino a 1/09/2021
# Append 9 new daily rows (2021-08-24 .. 2021-09-01) to the bottom of the
# single data sheet, continuing the existing A:D layout
# (data, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A carries the date-formatted style used by every existing row in
# column A (bordered, bold, centered, custom date number format). Copy that
# formatting from the last existing row (357) instead of re-building it,
# so the new cells land on the very same style index.
$formatSource = $ws.Range("A357")

$rows = @(
    @(358, 44432, 0, 12, 77.50435962022863),
    @(359, 44433, 0, 10, 64.58696635019054),
    @(360, 44434, 2, 12, 77.50435962022863),
    @(361, 44435, 1, 9,  58.12826971517148),
    @(362, 44436, 2, 6,  38.75217981011431),
    @(363, 44437, 0, 6,  38.75217981011431),
    @(364, 44438, 0, 5,  32.29348317509527),
    @(365, 44439, 0, 5,  32.29348317509527),
    @(366, 44440, 0, 5,  32.29348317509527)
)

foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    $formatSource.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
